# Update column F (dSF) values to match repulled data / recalculated mean values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1
$ws.Range("F3").Value = -4
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = 3
$ws.Range("F7").Value = -1
$ws.Range("F12").Value = -10
$ws.Range("F13").Value = -8
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 1
$ws.Range("F16").Value = -4
$ws.Range("F17").Value = -1
